# Apply the "MacroLibrary" workbook edit:
#  - Insert a new json-library entry "compact(var,json,removeEmpty)" at the
#    top of the `json` list (column M of the hidden "#system" sheet),
#    pushing the existing entries down by one row.
#  - Insert a new web-library entry "clickAll(locator)" into the `web` list
#    (column Y of the hidden "#system" sheet) right after "click(locator)",
#    pushing the existing entries down by one row.
#  - Extend the `json` and `web` defined names to cover the newly added rows.
#
# NOTE: these two lists live in different columns of the same worksheet, so a
# normal whole-row insert would incorrectly shift every other column's data
# too. Instead we manually shift just the affected column's cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) json list (column M = 13): currently M2:M16, becomes M2:M17.
#    New entry goes in at M12 ("compact(var,json,removeEmpty)" sorts
#    alphabetically between "beautify(json,var)" and
#    "fromCsv(csv,header,jsonFile)").
# ---------------------------------------------------------------------
$jsonCol = 13
for ($r = 17; $r -ge 13; $r--) {
    $srcVal = $ws.Cells.Item($r - 1, $jsonCol).Value()
    $ws.Cells.Item($r, $jsonCol).Value = $srcVal
}
$ws.Cells.Item(12, $jsonCol).Value = "compact(var,json,removeEmpty)"

# ---------------------------------------------------------------------
# 2) web list (column Y = 25): currently Y2:Y128, becomes Y2:Y129.
#    New entry goes in at Y50 ("clickAll(locator)" sorts alphabetically
#    between "click(locator)" and "clickAndWait(locator,waitMs)").
# ---------------------------------------------------------------------
$webCol = 25
for ($r = 129; $r -ge 51; $r--) {
    $srcVal = $ws.Cells.Item($r - 1, $webCol).Value()
    $ws.Cells.Item($r, $webCol).Value = $srcVal
}
$ws.Cells.Item(50, $webCol).Value = "clickAll(locator)"

# ---------------------------------------------------------------------
# 3) Update the defined names so they cover the extra row each.
# ---------------------------------------------------------------------
$jsonName = $wb.Names.Item("json")
$jsonName.RefersTo = "='#system'!`$M`$2:`$M`$17"

$webName = $wb.Names.Item("web")
$webName.RefersTo = "='#system'!`$Y`$2:`$Y`$129"
